$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.735.51"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "'3.150.13"
$ws.Range("E3").Value = "  +1.52%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'536.15"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "'143.81"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'3.148.85"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").Value = "'7.19"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("D13").Value = "'3.685.19"
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("E14").Value = "  +3.25%  "

$ws.Range("D15").Value = "'25.88"
$ws.Range("E15").Value = "  -3.99%  "

$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").Value = "'58.759.13"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").Value = "'3.145.74"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "'12.95"
$ws.Range("E20").Value = "  +0.26%  "

$ws.Range("D21").Value = "'8.02"
$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("D22").Value = "'343.68"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  +1.80%  "

$ws.Range("D25").Value = "'67.94"
$ws.Range("E25").Value = "  +2.72%  "

$ws.Range("E26").Value = "  -0.23%  "

$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'0.0₃0942"
$ws.Range("E28").Value = "  +2.59%  "

$ws.Range("D29").Value = "'7.59"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("D30").Value = "'6.51"
$ws.Range("E30").Value = "  -2.06%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  +1.76%  "

$ws.Range("D33").Value = "'21.20"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("D35").Value = "'4.83"
$ws.Range("E35").Value = "  +3.05%  "

$ws.Range("D36").Value = "'158.08"
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("E37").Value = "  +3.43%  "

$ws.Range("D38").Value = "'26.30"
$ws.Range("E38").Value = "  -2.41%  "

$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").Value = "'1.68"
$ws.Range("E40").Value = "  +12.32%  "

$ws.Range("E41").Value = "  -0.62%  "

$ws.Range("D42").Value = "'0.710"
$ws.Range("E42").Value = "  +4.55%  "

$ws.Range("E43").Value = "  +3.55%  "

$ws.Range("D44").Value = "'3.184.02"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D45").Value = "'36.85"
$ws.Range("E45").Value = "  +0.08%  "

$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  +3.37%  "

$ws.Range("D48").Value = "'2.315.07"
$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("D49").Value = "'1.02"
$ws.Range("E49").Value = "  +5.15%  "

$ws.Range("D50").Value = "'20.82"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").Value = "'6.10"
$ws.Range("E51").Value = "  +1.83%  "
